# Change store locator function code in response to changes on Sprouts' site
#
# Semantic changes applied:
#  - Remove the scraped price (C2) and image-link URL (D2) values that were
#    populated by the old store-locator code, while keeping C2's number
#    formatting/style intact.
#  - Remove the picture that had been inserted next to the product row.
#  - Update the current selection on the "Sprouts" sheet from F2 to D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprouts")

# Drop the picture that was anchored near F2 (no longer produced by the
# updated locator code).
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete() | Out-Null
}

# Clear the scraped price and image URL cells (C2, D2) while preserving
# the existing cell formatting on C2.
$ws.Range("C2").ClearContents() | Out-Null
$ws.Range("D2").ClearContents() | Out-Null

# Update the saved selection/active cell.
$null = $ws.Range("D4").Select()
